# Update heterogeneity meta-analysis table: add PdeltaAIC as a covariate
# for the CG path. This changes several Q (heterogeneity statistic) and
# pval values throughout the "Temperature" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A = Q, Column B = I2, Column F = pval (all stored as text values,
# including leading-space padding). Force Text number format first so
# Excel does not auto-convert these into numbers and strip formatting.

$cells = @("A2","A3","B3","A4","B4","A5","B5","A6","F6","A7","B7","F7","A10","B10","A12","B12","F12","A13")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A2").Value = "1618.4"

$ws.Range("A3").Value = " 249.0"
$ws.Range("B3").Value = "0.637"

$ws.Range("A4").Value = " 251.9"
$ws.Range("B4").Value = "0.545"

$ws.Range("A5").Value = " 345.8"
$ws.Range("B5").Value = "0.730"

$ws.Range("A6").Value = " 124.0"
$ws.Range("F6").Value = "0.078"

$ws.Range("A7").Value = " 121.2"
$ws.Range("B7").Value = "0.286"
$ws.Range("F7").Value = "0.107"

$ws.Range("A10").Value = " 121.9"
$ws.Range("B10").Value = "0.390"

$ws.Range("A12").Value = "  89.1"
$ws.Range("B12").Value = "0.062"
$ws.Range("F12").Value = "0.506"

$ws.Range("A13").Value = "  55.6"
